# M11 Froze Token Embeddings
# Update the B (predicted token) and C (score) columns in the ASR Results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "<sabet>"
$ws.Range("C2").Value = 38

$ws.Range("C3").Value = 31

$ws.Range("B4").Value = "<mind>"
$ws.Range("C4").Value = 29

$ws.Range("C5").Value = 36

$ws.Range("B6").Value = "<are>"
$ws.Range("C6").Value = 34

$ws.Range("C7").Value = 31

$ws.Range("B8").Value = "<thow>"
$ws.Range("C8").Value = 25

$ws.Range("C9").Value = 28

$ws.Range("B10").Value = "<time>"
$ws.Range("C10").Value = 30

$ws.Range("C11").Value = 33

$ws.Range("B12").Value = "<like>"
$ws.Range("C12").Value = 27

$ws.Range("B13").Value = "<can>"
$ws.Range("C13").Value = 32

$ws.Range("C14").Value = 31

$ws.Range("B15").Value = "<up>"
$ws.Range("C15").Value = 32

$ws.Range("B16").Value = "<towe>"
$ws.Range("C16").Value = 17
